$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.151.17'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '3.395.05'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '566.65'
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").Value = '155.53'
$ws.Range("E6").Value = '  +1.94%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.396.28'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +2.67%  '
$ws.Range("D10").Value = '7.41'
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").Value = '3.984.35'
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("E14").Value = '  -3.17%  '
$ws.Range("D15").Value = '0.0000190'
$ws.Range("E15").Value = '  +5.85%  '
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = '63.227.53'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '3.433.44'
$ws.Range("E18").Value = '  +2.40%  '
$ws.Range("E19").Value = '  -1.93%  '
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '377.88'
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("E22").Value = '  -3.44%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '71.37'
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").Value = '0.528'
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  +24.47%  '
$ws.Range("D27").Value = '9.45'
$ws.Range("E27").Value = '  +6.51%  '
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '6.07'
$ws.Range("E30").Value = '  +7.68%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '1.36'
$ws.Range("E31").Value = '  +4.12%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '2.00'
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("D33").Value = '23.15'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = '6.79'
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").Value = '159.90'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").Value = '2.951.86'
$ws.Range("E39").Value = '  +4.75%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '27.00'
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.0753'
$ws.Range("E41").Value = '  +1.23%  '
$ws.Range("D42").Value = '1.82'
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("D45").Value = '0.758'
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("D47").Value = '23.31'
$ws.Range("E47").Value = '  +6.34%  '
$ws.Range("E48").Value = '  +3.54%  '
$ws.Range("E49").Value = '  +20.38%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").Value = '0.830'
$ws.Range("E51").Value = '  +3.90%  '
